# Remove the two trailing exercise sections ("Случаен списък" / random list,
# and "Поредица от стрингове" / string stack) from the Inheritance Basics
# exercises document, together with their explanatory paragraphs and the two
# screenshots that illustrate them. Everything from the "Случаен списък"
# Heading2 paragraph through to the end of the document body is removed; the
# paragraph right before it (the screenshot that closes out the previous
# "Поредица от животни" section) is left untouched, and the section
# properties (page size/margins/header/footer) immediately follow it.

$d = $word.ActiveDocument

# Anchor on the heading text that starts the block we need to drop - this is
# more robust than hard-coding paragraph indices.
$finder = $d.Content
$found = $finder.Find.Execute("Случаен списък", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Случаен списък' heading to anchor the deletion."
}

# Expand back to the very start of that heading paragraph, then extend all
# the way to the end of the document body so every following paragraph
# (including the final screenshot) is captured too.
$deleteStart = $finder.Paragraphs.Item(1).Range.Start
$deleteEnd = $d.Content.End

$target = $d.Range($deleteStart, $deleteEnd)
$target.Delete()

Write-Output "Remaining paragraphs: $($d.Paragraphs.Count)"
